$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 16.458
$ws.Range("A3").Value = -21.879
$ws.Range("A14").Value = -21.659
$ws.Range("A16").Value = -22.051
$ws.Range("E18").Value = 16.481
$ws.Range("A21").Value = -19.947
$ws.Range("A23").Value = -20.302
$ws.Range("E24").Value = 16.963
$ws.Range("A25").Value = -21.632
$ws.Range("E25").Value = 17.259
$ws.Range("A26").Value = -21.089
$ws.Range("E27").Value = 16.55
$ws.Range("A29").Value = -21.04
$ws.Range("E30").Value = 16.466
$ws.Range("E31").Value = 16.605
$ws.Range("E39").Value = 16.552
$ws.Range("A40").Value = -20.174
$ws.Range("E42").Value = 16.782
$ws.Range("E48").Value = 17.179
$ws.Range("E51").Value = 16.617
$ws.Range("E52").Value = 16.543
$ws.Range("A53").Value = -21.692
$ws.Range("E55").Value = 16.413
$ws.Range("E56").Value = 16.303
$ws.Range("A57").Value = -22.213
$ws.Range("E57").Value = 16.453
$ws.Range("A59").Value = -22.5
$ws.Range("E60").Value = 16.592
$ws.Range("A65").Value = -21.533
$ws.Range("A69").Value = -21.601
$ws.Range("E73").Value = 16.572
$ws.Range("E74").Value = 16.653
$ws.Range("A79").Value = -21.167
$ws.Range("A83").Value = -22.035
$ws.Range("E89").Value = 17.362
$ws.Range("E90").Value = 16.23
$ws.Range("A91").Value = -21.533
$ws.Range("E92").Value = 17.493
$ws.Range("A93").Value = -21.22
$ws.Range("A100").Value = -21.955
